$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand Table1 by one row (extends table ref + autoFilter ref A1:C27 -> A1:C28)
$tbl = $ws.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null

# Populate the new row's cells (row 28)
$ws.Range("A28").Value2 = "Distributed"
$ws.Range("B28").Value2 = "Distributed"
$ws.Range("C28").Value2 = "Both horizontally and vertically distributed."

# Apply the new "distributed" horizontal + vertical alignment (with wrap) to C28
$c = $ws.Range("C28")
$c.HorizontalAlignment = -4117   # xlDistributed
$c.VerticalAlignment = -4117     # xlDistributed
$c.WrapText = $true

# Match the custom row height used for the new row
$ws.Rows.Item(28).RowHeight = 43.5

# Update the view: scroll down a bit and select the new last cell
$ws.Range("C28").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 19
